# Generate Report for Handback
# Adds a new handback row (file 7d65ac3a-205f-4b2c-83a2-ea551215c72e.md) to the
# Overview sheet plus the per-locale zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$fileId   = "7d65ac3a-205f-4b2c-83a2-ea551215c72e"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$status   = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276  # decimal OLE BGR for FF6495ED (matches existing HyperLink font color)

function Style-AsHyperlink($rng) {
    $rng.Style = "Hyperlink"
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
}

function Style-AsDateText($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# The source workbook stores "True"/"False"/"" as literal shared-string text
# (not native booleans). A leading apostrophe forces Excel's COM layer to
# commit a text value instead of auto-coercing "True"/"False" to booleans
# (or dropping an all-blank assignment entirely); resetting the style back to
# Normal afterwards clears the quote-prefix marker so the cell lands exactly
# like a plain, unstyled text cell (matching the source formatting).
function Set-TextValue($rng, $text) {
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3, 1).Value = $mdName                                  # File Name
$wsOverview.Cells.Item(3, 2).Value = $mdPath                                  # Path And Name
$wsOverview.Cells.Item(3, 3).Value = ".md"                                    # Extension
# Publish URL (column 4) left blank, matches source row
$wsOverview.Cells.Item(3, 5).Value = $status                                  # zh-cn
$wsOverview.Cells.Item(3, 6).Value = $status                                  # de-de
$wsOverview.Cells.Item(3, 7).Value = "2016-10-27 09:15:08"                    # Latest HO Xliff Generate Date
Style-AsDateText $wsOverview.Cells.Item(3, 7)

Style-AsHyperlink $wsOverview.Cells.Item(3, 2)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(3, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1563339901e18b51bb252728e9ba7d50e9cf63/$mdPath",
    "",
    "",
    $mdPath
) | Out-Null

# ---------------------------------------------------------------------------
# Shared per-locale row writer (zh-cn / de-de detail sheets)
# ---------------------------------------------------------------------------
function Add-LocaleHandbackRow($sheetName, $repoSuffix, $handoffFile, $handoffDate, $handbackDate) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null
    $r = 3

    $ws.Cells.Item($r, 1).Value  = $mdName            # Source File Name
    $ws.Cells.Item($r, 2).Value  = ".md"              # File Extension
    $ws.Cells.Item($r, 3).Value  = $status            # Status
    $ws.Cells.Item($r, 4).Value  = "e2e"              # Source Path
    $ws.Cells.Item($r, 5).Value  = "ht"               # Priority
    Set-TextValue $ws.Cells.Item($r, 6) "True"        # Content Duplicate
    $ws.Cells.Item($r, 7).Value  = $handoffFile       # Correspond Handoff File
    $ws.Cells.Item($r, 8).Value  = $handoffDate       # Correspond Handoff Datetime
    $ws.Cells.Item($r, 9).Value  = $mdName            # Target File
    $ws.Cells.Item($r, 10).Value = $handoffFile       # Correspond Handback File
    $ws.Cells.Item($r, 11).Value = $handbackDate      # Correspond Handback DateTime
    Set-TextValue $ws.Cells.Item($r, 12) ""           # Reference Tokens
    Set-TextValue $ws.Cells.Item($r, 13) "True"       # To be localized
    Set-TextValue $ws.Cells.Item($r, 14) ""           # Dependency From
    Set-TextValue $ws.Cells.Item($r, 15) "False"      # Has metadata
    Set-TextValue $ws.Cells.Item($r, 16) ""           # Error Detail

    Style-AsDateText $ws.Cells.Item($r, 8)
    Style-AsDateText $ws.Cells.Item($r, 11)

    Style-AsHyperlink $ws.Cells.Item($r, 1)
    Style-AsHyperlink $ws.Cells.Item($r, 9)

    $ws.Hyperlinks.Add(
        $ws.Cells.Item($r, 1),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1563339901e18b51bb252728e9ba7d50e9cf63/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Cells.Item($r, 9),
        "https://github.com/OpenLocalizationTestOrg/ol-test0-$repoSuffix/blob/f49149d5917f70e585a7286f1e15035f3021364b/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
Add-LocaleHandbackRow "zh-cn" "zhcn" `
    "$fileId.714ebb150d80af156cbfa7cc82ee75ab27f3a527.zh-cn.xlf" `
    "2016-10-27 09:14:55" `
    "2016-10-27 09:15:36"

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
Add-LocaleHandbackRow "de-de" "dede" `
    "$fileId.714ebb150d80af156cbfa7cc82ee75ab27f3a527.de-de.xlf" `
    "2016-10-27 09:15:08" `
    "2016-10-27 09:15:52"
